{"js": "// Office.js (Word JavaScript API) script.\n// Adds a new \"4) Evaluate Each Potential Solution\" section (with a-d\n// sub points) plus a new \"5) Choose a solution...\" line right after the\n// existing \"Identify Potential Solutions\" list for Problem #1 (the\n// cat/parrot/seed river-crossing problem), before the \"Socks in the\n// Dark\" section.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph: \"d) Try to find a bigger boat or get\n// help from someone.\" -- the last bullet of the \"3) Identify Potential\n// Solutions\" list for Problem #1.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"Try to find a bigger boat or get help from someone.\"\n    ) !== -1\n  ) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert the new paragraphs, in order, right after the anchor\n// paragraph.\nlet cursor = anchor;\ncursor = cursor.insertParagraph(\"\", Word.InsertLocation.after);\ncursor = cursor.insertParagraph(\n  \"4) Evaluate Each Potential Solution\",\n  Word.InsertLocation.after\n);\n\nconst lettered = [\n  \"a) This solution would work as the man would be able to eliminate the constraint of worrying about the safety of his cat, bird, and bag of seed while transporting them.\",\n  \"b) This would be the optimal solution, provided that there is actually a way to get everything across the river safely without any additional equipment or help.\",\n  \"c) This solution is not optimal as the man would lose one of his things during the transporting and the goal is to get all three things across safely.\",\n  \"d) This solution has potential, but it is working outside of the context of the problem itself. I\u2019m sure that if there were a source of help readily available the man would not be faced with this problem to begin with.\",\n];\n\nfor (const text of lettered) {\n  cursor = cursor.insertParagraph(text, Word.InsertLocation.after);\n  cursor.leftIndent = 36; // 0.5in == 720 twips, matches the existing lettered lists\n  cursor.firstLineIndent = 0;\n}\n\nawait context.sync();\n\n// The pre-existing blank paragraph that used to sit directly before\n// \"Socks in the Dark\" becomes the new \"5) Choose a solution...\" line.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet fiveIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (\n    paragraphs2.items[i].text.indexOf(\n      \"d) This solution has potential, but it is working outside\"\n    ) !== -1\n  ) {\n    fiveIndex = i + 1;\n    break;\n  }\n}\n\nif (fiveIndex === -1) {\n  throw new Error(\"Paragraph before 'Socks in the Dark' not found\");\n}\n\nconst fiveParagraph = paragraphs2.items[fiveIndex];\nfiveParagraph.insertText(\n  \"5) Choose a solution and develop a plan to implement it.\",\n  Word.InsertLocation.replace\n);\nfiveParagraph.leftIndent = 0;\nfiveParagraph.firstLineIndent = 0;\n\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) script.\n# Adds a new \"4) Evaluate Each Potential Solution\" section (with a-d sub\n# points) plus a new \"5) Choose a solution...\" line right after the\n# existing \"Identify Potential Solutions\" list for Problem #1 (the\n# cat/parrot/seed river-crossing problem), before the \"Socks in the Dark\"\n# section.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: \"d) Try to find a bigger boat or get help\n# from someone.\" -- the last bullet of the \"3) Identify Potential\n# Solutions\" list for Problem #1.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*Try to find a bigger boat or get help from someone.*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Insert the new paragraphs, in order, right after the anchor paragraph.\n$newTexts = @(\n    \"\",\n    \"4) Evaluate Each Potential Solution\",\n    \"a) This solution would work as the man would be able to eliminate the constraint of worrying about the safety of his cat, bird, and bag of seed while transporting them.\",\n    \"b) This would be the optimal solution, provided that there is actually a way to get everything across the river safely without any additional equipment or help.\",\n    \"c) This solution is not optimal as the man would lose one of his things during the transporting and the goal is to get all three things across safely.\",\n    \"d) This solution has potential, but it is working outside of the context of the problem itself. I\u2019m sure that if there were a source of help readily available the man would not be faced with this problem to begin with.\"\n)\n\n$insertAfterIndex = $anchorIndex\nforeach ($t in $newTexts) {\n    $p = $d.Paragraphs($insertAfterIndex)\n    $p.Range.InsertParagraphAfter()\n    $insertAfterIndex = $insertAfterIndex + 1\n    if ($t -ne \"\") {\n        $newP = $d.Paragraphs($insertAfterIndex)\n        $newP.Range.Text = $t\n    }\n}\n\n# The four lettered evaluations (a-d) get a 0.5in left indent with no\n# first-line indent, matching the existing lettered list style elsewhere\n# in the document.\nfor ($i = $anchorIndex + 3; $i -le $anchorIndex + 6; $i++) {\n    $p = $d.Paragraphs($i)\n    $p.Format.LeftIndent = 36\n    $p.Format.FirstLineIndent = 0\n}\n\n# The pre-existing blank paragraph that used to sit directly before\n# \"Socks in the Dark\" becomes the new \"5) Choose a solution...\" line.\n$fiveIndex = $insertAfterIndex + 1\n$fiveP = $d.Paragraphs($fiveIndex)\n$fiveP.Range.Text = \"5) Choose a solution and develop a plan to implement it.\"\n$fiveP.Format.LeftIndent = 0\n$fiveP.Format.FirstLineIndent = 0\n"}
